$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02; $ws.Range("C2").Value = 1.014472626724884; $ws.Range("D2").Value = 1.04003660180095; $ws.Range("E2").Value = 1.016200295271653; $ws.Range("F2").Value = 1.042866889444925; $ws.Range("I2").Value = 1.033320481670482; $ws.Range("J2").Value = 1.019703665611906; $ws.Range("K2").Value = 1.042820089625907; $ws.Range("L2").Value = 1.019053176183402; $ws.Range("M2").Value = 1.045642366384276; $ws.Range("N2").Value = 1.010512101947724
$ws.Range("B3").Value = 1.02; $ws.Range("C3").Value = 1.015984545280922; $ws.Range("D3").Value = 1.040703179401188; $ws.Range("E3").Value = 1.017501368530942; $ws.Range("F3").Value = 1.043884044043582; $ws.Range("I3").Value = 1.033432032108268; $ws.Range("J3").Value = 1.020847203927548; $ws.Range("K3").Value = 1.043297387072164; $ws.Range("L3").Value = 1.020158205455324; $ws.Range("M3").Value = 1.046469899985069; $ws.Range("N3").Value = 1.010906688163961
$ws.Range("B4").Value = 1.02; $ws.Range("C4").Value = 1.016961764608468; $ws.Range("D4").Value = 1.04113361854109; $ws.Range("E4").Value = 1.018342528210695; $ws.Range("F4").Value = 1.044541364715324; $ws.Range("I4").Value = 1.033502509462256; $ws.Range("J4").Value = 1.021585755493524; $ws.Range("K4").Value = 1.043604669792253; $ws.Range("L4").Value = 1.020871982589181; $ws.Range("M4").Value = 1.047003919362935; $ws.Range("N4").Value = 1.011161095301501
$ws.Range("B5").Value = 1.02; $ws.Range("C5").Value = 1.01737233355766; $ws.Range("D5").Value = 1.041314363960756; $ws.Range("E5").Value = 1.018695984428887; $ws.Range("F5").Value = 1.044817501088867; $ws.Range("I5").Value = 1.033531730264007; $ws.Range("J5").Value = 1.021895914100059; $ws.Range("K5").Value = 1.043733477897608; $ws.Range("L5").Value = 1.021171760040697; $ws.Range("M5").Value = 1.047228075233636; $ws.Range("N5").Value = 1.011267830293866
$ws.Range("B6").Value = 1.02; $ws.Range("C6").Value = 1.017441255213879; $ws.Range("D6").Value = 1.04134469951992; $ws.Range("E6").Value = 1.018755321630427; $ws.Range("F6").Value = 1.044863853840272; $ws.Range("I6").Value = 1.03353661264346; $ws.Range("J6").Value = 1.021947971986615; $ws.Range("K6").Value = 1.043755083441887; $ws.Range("L6").Value = 1.021222076870625; $ws.Range("M6").Value = 1.047265691745848; $ws.Range("N6").Value = 1.011285738846469
$ws.Range("B7").Value = 1.02; $ws.Range("C7").Value = 1.016967251639693; $ws.Range("D7").Value = 1.041136034500455; $ws.Range("E7").Value = 1.018347251760641; $ws.Range("F7").Value = 1.044545055252126; $ws.Range("I7").Value = 1.03350290151489; $ws.Range("J7").Value = 1.021589901132135; $ws.Range("K7").Value = 1.043606392400223; $ws.Range("L7").Value = 1.020875989381045; $ws.Range("M7").Value = 1.047006915902149; $ws.Range("N7").Value = 1.011162522353339
$ws.Range("B8").Value = 1.02; $ws.Range("C8").Value = 1.014983817610201; $ws.Range("D8").Value = 1.040262057123983; $ws.Range("E8").Value = 1.016640151808293; $ws.Range("F8").Value = 1.043210817690324; $ws.Range("I8").Value = 1.033358533260207; $ws.Range("J8").Value = 1.020090421614083; $ws.Range("K8").Value = 1.04298171772145; $ws.Range("L8").Value = 1.019426887959555; $ws.Range("M8").Value = 1.045922336303236; $ws.Range("N8").Value = 1.010645644977287
$ws.Range("B9").Value = 1.02; $ws.Range("C9").Value = 1.011480039302787; $ws.Range("D9").Value = 1.038715258406494; $ws.Range("E9").Value = 1.013626250248608; $ws.Range("F9").Value = 1.040853174920126; $ws.Range("I9").Value = 1.033091097909824; $ws.Range("J9").Value = 1.017437231075015; $ws.Range("K9").Value = 1.041869004029306; $ws.Range("L9").Value = 1.016863584155789; $ws.Range("M9").Value = 1.044000001871817; $ws.Range("N9").Value = 1.009727747090801
$ws.Range("B10").Value = 1.02; $ws.Range("C10").Value = 1.009137847961366; $ws.Range("D10").Value = 1.037679537571036; $ws.Range("E10").Value = 1.011612755257051; $ws.Range("F10").Value = 1.03927691911508; $ws.Range("I10").Value = 1.032904045206979; $ws.Range("J10").Value = 1.015660763196928; $ws.Range("K10").Value = 1.041119146068815; $ws.Range("L10").Value = 1.015147810940389; $ws.Range("M10").Value = 1.042710850355354; $ws.Range("N10").Value = 1.009110940538285
$ws.Range("B11").Value = 1.02; $ws.Range("C11").Value = 1.008122042420112; $ws.Range("D11").Value = 1.037229987276775; $ws.Range("E11").Value = 1.010739809565534; $ws.Range("F11").Value = 1.038593295075287; $ws.Range("I11").Value = 1.032820971533971; $ws.Range("J11").Value = 1.014889637731391; $ws.Range("K11").Value = 1.040792537045736; $ws.Range("L11").Value = 1.014403154310602; $ws.Range("M11").Value = 1.042150811403045; $ws.Range("N11").Value = 1.008842676582181
$ws.Range("B12").Value = 1.02; $ws.Range("C12").Value = 1.007744474137259; $ws.Range("D12").Value = 1.037062842502643; $ws.Range("E12").Value = 1.010415388300099; $ws.Range("F12").Value = 1.038339200145285; $ws.Range("I12").Value = 1.032789801959818; $ws.Range("J12").Value = 1.014602914931416; $ws.Range("K12").Value = 1.040670931726069; $ws.Range("L12").Value = 1.014126291813578; $ws.Range("M12").Value = 1.041942511478225; $ws.Range("N12").Value = 1.008742851637598
$ws.Range("B13").Value = 1.02; $ws.Range("C13").Value = 1.007825475408861; $ws.Range("D13").Value = 1.037098702946274; $ws.Range("E13").Value = 1.010484985581159; $ws.Range("F13").Value = 1.038393711938773; $ws.Range("I13").Value = 1.032796502072577; $ws.Range("J13").Value = 1.01466443128562; $ws.Range("K13").Value = 1.040697029527468; $ws.Range("L13").Value = 1.014185691798605; $ws.Range("M13").Value = 1.041987205063376; $ws.Range("N13").Value = 1.008764272595604
$ws.Range("B14").Value = 1.02; $ws.Range("C14").Value = 1.008090837720714; $ws.Range("D14").Value = 1.037216174338332; $ws.Range("E14").Value = 1.010712996317184; $ws.Range("F14").Value = 1.038572294917992; $ws.Range("I14").Value = 1.03281840141851; $ws.Range("J14").Value = 1.014865943139562; $ws.Range("K14").Value = 1.04078249099333; $ws.Range("L14").Value = 1.014380274188718; $ws.Range("M14").Value = 1.042133598918406; $ws.Range("N14").Value = 1.008834428700447
$ws.Range("B15").Value = 1.02; $ws.Range("C15").Value = 1.008254302411132; $ws.Range("D15").Value = 1.037288530900274; $ws.Range("E15").Value = 1.010853458478893; $ws.Range("F15").Value = 1.038682303658233; $ws.Range("I15").Value = 1.032831852942929; $ws.Range("J15").Value = 1.0149900622922; $ws.Range("K15").Value = 1.040835108423488; $ws.Range("L15").Value = 1.014500127671839; $ws.Range("M15").Value = 1.04222376030457; $ws.Range("N15").Value = 1.008877630308489
$ws.Range("B16").Value = 1.02; $ws.Range("C16").Value = 1.009205228679451; $ws.Range("D16").Value = 1.037709350060756; $ws.Range("E16").Value = 1.011670666267362; $ws.Range("F16").Value = 1.039322265750767; $ws.Range("I16").Value = 1.032909514737067; $ws.Range("J16").Value = 1.015711899600034; $ws.Range("K16").Value = 1.041140781641946; $ws.Range("L16").Value = 1.015197194682216; $ws.Range("M16").Value = 1.042747979623351; $ws.Range("N16").Value = 1.009128719238746
$ws.Range("B17").Value = 1.02; $ws.Range("C17").Value = 1.00980127900544; $ws.Range("D17").Value = 1.037973030592478; $ws.Range("E17").Value = 1.012182982799052; $ws.Range("F17").Value = 1.039723402510416; $ws.Range("I17").Value = 1.03295767340871; $ws.Range("J17").Value = 1.016164174942012; $ws.Range("K17").Value = 1.041332009250274; $ws.Range("L17").Value = 1.015633982795791; $ws.Range("M17").Value = 1.043076317848331; $ws.Range("N17").Value = 1.009285902594273
$ws.Range("B18").Value = 1.02; $ws.Range("C18").Value = 1.01014878928281; $ws.Range("D18").Value = 1.038126727141745; $ws.Range("E18").Value = 1.012481703743412; $ws.Range("F18").Value = 1.03995727329621; $ws.Range("I18").Value = 1.03298556299309; $ws.Range("J18").Value = 1.016427796315401; $ws.Range("K18").Value = 1.041443364281901; $ws.Range("L18").Value = 1.015888588808388; $ws.Range("M18").Value = 1.04326765575171; $ws.Range("N18").Value = 1.009377471019994
$ws.Range("B19").Value = 1.02; $ws.Range("C19").Value = 1.010267255161924; $ws.Range("D19").Value = 1.038179116080815; $ws.Range("E19").Value = 1.01258354232336; $ws.Range("F19").Value = 1.040036999356334; $ws.Range("I19").Value = 1.032995038603634; $ws.Range("J19").Value = 1.016517653517244; $ws.Range("K19").Value = 1.041481302136448; $ws.Range("L19").Value = 1.015975375059163; $ws.Range("M19").Value = 1.043332867197682; $ws.Range("N19").Value = 1.009408674218924
$ws.Range("B20").Value = 1.02; $ws.Range("C20").Value = 1.0097373446283; $ws.Range("D20").Value = 1.037944750908557; $ws.Range("E20").Value = 1.012128026954085; $ws.Range("F20").Value = 1.039680375259582; $ws.Range("I20").Value = 1.032952527183476; $ws.Range("J20").Value = 1.016115669072883; $ws.Range("K20").Value = 1.041311511455017; $ws.Range("L20").Value = 1.015587136686913; $ws.Range("M20").Value = 1.04304110850622; $ws.Range("N20").Value = 1.009269050111003
$ws.Range("B21").Value = 1.02; $ws.Range("C21").Value = 1.008012702189515; $ws.Range("D21").Value = 1.03718158640783; $ws.Range("E21").Value = 1.010645857595856; $ws.Range("F21").Value = 1.038519711294996; $ws.Range("I21").Value = 1.032811961228563; $ws.Range("J21").Value = 1.014806611047435; $ws.Range("K21").Value = 1.040757332680301; $ws.Range("L21").Value = 1.014322981858226; $ws.Range("M21").Value = 1.042090497225722; $ws.Range("N21").Value = 1.008813774455461
$ws.Range("B22").Value = 1.02; $ws.Range("C22").Value = 1.006926883318906; $ws.Range("D22").Value = 1.036700818543317; $ws.Range("E22").Value = 1.00971296982019; $ws.Range("F22").Value = 1.037788991312763; $ws.Range("I22").Value = 1.032721774868431; $ws.Range("J22").Value = 1.013981858907653; $ws.Range("K22").Value = 1.040407230924041; $ws.Range("L22").Value = 1.01352662780813; $ws.Range("M22").Value = 1.041491209344009; $ws.Range("N22").Value = 1.008526483458712
$ws.Range("B23").Value = 1.02; $ws.Range("C23").Value = 1.007502638691648; $ws.Range("D23").Value = 1.036955771440791; $ws.Range("E23").Value = 1.010207607298348; $ws.Range("F23").Value = 1.038176451908121; $ws.Range("I23").Value = 1.032769755653644; $ws.Range("J23").Value = 1.014419238521985; $ws.Range("K23").Value = 1.040592984690715; $ws.Range("L23").Value = 1.013948937217806; $ws.Range("M23").Value = 1.041809055574957; $ws.Range("N23").Value = 1.008678881239402
$ws.Range("B24").Value = 1.02; $ws.Range("C24").Value = 1.009766234302906; $ws.Range("D24").Value = 1.037957529601723; $ws.Range("E24").Value = 1.012152859458416; $ws.Range("F24").Value = 1.039699817748601; $ws.Range("I24").Value = 1.032954853160979; $ws.Range("J24").Value = 1.016137587353471; $ws.Range("K24").Value = 1.041320774097701; $ws.Range("L24").Value = 1.015608304937432; $ws.Range("M24").Value = 1.043057018637279; $ws.Range("N24").Value = 1.009276665375042
$ws.Range("B25").Value = 1.02; $ws.Range("C25").Value = 1.012386931909205; $ws.Range("D25").Value = 1.039115940522589; $ws.Range("E25").Value = 1.014406135652801; $ws.Range("F25").Value = 1.041463467260606; $ws.Range("I25").Value = 1.033161781135759; $ws.Range("J25").Value = 1.018124472367177; $ws.Range("K25").Value = 1.042158085083129; $ws.Range("L25").Value = 1.017527453380116; $ws.Range("M25").Value = 1.044498303604836; $ws.Range("N25").Value = 1.009965897125824

Write-Host "applied 264 cell updates"
